$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5165
$ws.Range("J18").Value = 24998.5
$ws.Range("L18").Value = 24998.5
$ws.Range("N18").Value = -25566.5
$ws.Range("H40").Value = 10749.714
$ws.Range("I40").Value = 11917.818
$ws.Range("J40").Value = 6466.6665
$ws.Range("K40").Value = 11917.818
$ws.Range("L40").Value = 6466.6665
$ws.Range("M40").Value = -11742.818
$ws.Range("N40").Value = -6816.6665
$ws.Range("H43").Value = 7599.75
$ws.Range("J43").Value = 2579.6
$ws.Range("L43").Value = 2579.6
$ws.Range("N43").Value = -2717.6
$ws.Range("H64").Value = 5749.5
$ws.Range("I64").Value = 4831.6665
$ws.Range("J64").Value = 6142.857
$ws.Range("K64").Value = 4831.6665
$ws.Range("L64").Value = 6142.857
$ws.Range("M64").Value = -4583.6665
$ws.Range("N64").Value = -6638.857
$ws.Range("H67").Value = 5749.5
$ws.Range("I67").Value = 4831.6665
$ws.Range("J67").Value = 6142.857
$ws.Range("K67").Value = 4831.6665
$ws.Range("L67").Value = 6142.857
$ws.Range("M67").Value = -3973.6665
$ws.Range("N67").Value = -7858.857
$ws.Range("I74").Value = 125006504
$ws.Range("K74").Value = 125006504
$ws.Range("M74").Value = -125005568
$ws.Range("I77").Value = 125006504
$ws.Range("K77").Value = 625032520
$ws.Range("M77").Value = -625027840
$ws.Range("H106").Value = 2506.889
$ws.Range("I106").Value = 2392.9412
$ws.Range("K106").Value = 2392.9412
$ws.Range("M106").Value = -1761.9412
$ws.Range("H112").Value = 13266.267
$ws.Range("J112").Value = 13266.267
$ws.Range("L112").Value = 39798.801
$ws.Range("N112").Value = -42014.801
$ws.Range("H132").Value = 1972.4117
$ws.Range("I132").Value = 2007.0303
$ws.Range("J132").Value = 830
$ws.Range("K132").Value = 6021.090899999999
$ws.Range("L132").Value = 2490
$ws.Range("M132").Value = -3491.090899999999
$ws.Range("N132").Value = -7550
$ws.Range("H137").Value = 3245.0571
$ws.Range("I137").Value = 3287.524
$ws.Range("J137").Value = 3181.3572
$ws.Range("K137").Value = 9862.572
$ws.Range("L137").Value = 9544.071599999999
$ws.Range("M137").Value = -7312.572
$ws.Range("N137").Value = -14644.0716
$ws.Range("H138").Value = 1856694
$ws.Range("I138").Value = 4824.9
$ws.Range("J138").Value = 2277573.2
$ws.Range("K138").Value = 14474.7
$ws.Range("L138").Value = 6832719.600000001
$ws.Range("M138").Value = -9334.699999999999
$ws.Range("N138").Value = -6842999.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3041.5789
$ws.Range("I2").Value = 2893
$ws.Range("K2").Value = 2893
$ws.Range("M2").Value = -2780
$ws.Range("H32").Value = 1278035.9
$ws.Range("I32").Value = 1472599.4
$ws.Range("J32").Value = 5890.3076
$ws.Range("K32").Value = 1472599.4
$ws.Range("L32").Value = 5890.3076
$ws.Range("M32").Value = -1472312.4
$ws.Range("N32").Value = -6464.3076
$ws.Range("H45").Value = 3171.0833
$ws.Range("I45").Value = 2969.875
$ws.Range("J45").Value = 3573.5
$ws.Range("K45").Value = 2969.875
$ws.Range("L45").Value = 3573.5
$ws.Range("M45").Value = -2592.875
$ws.Range("N45").Value = -4327.5
$ws.Range("H52").Value = 59373.668
$ws.Range("J52").Value = 59373.668
$ws.Range("L52").Value = 59373.668
$ws.Range("N52").Value = -60009.668
$ws.Range("H61").Value = 6268.255
$ws.Range("I61").Value = 3324
$ws.Range("J61").Value = 10831.85
$ws.Range("K61").Value = 3324
$ws.Range("L61").Value = 10831.85
$ws.Range("M61").Value = -3112
$ws.Range("N61").Value = -11255.85
$ws.Range("H74").Value = 24140.639
$ws.Range("I74").Value = 30105.086
$ws.Range("J74").Value = 6744.3335
$ws.Range("K74").Value = 30105.086
$ws.Range("L74").Value = 6744.3335
$ws.Range("M74").Value = -29231.086
$ws.Range("N74").Value = -8492.333500000001
$ws.Range("H77").Value = 24140.639
$ws.Range("I77").Value = 30105.086
$ws.Range("J77").Value = 6744.3335
$ws.Range("K77").Value = 150525.43
$ws.Range("L77").Value = 33721.6675
$ws.Range("M77").Value = -146157.43
$ws.Range("N77").Value = -42457.6675
$ws.Range("H116").Value = 3041.5789
$ws.Range("I116").Value = 2893
$ws.Range("K116").Value = 2893
$ws.Range("M116").Value = -599
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080
$ws.Range("H136").Value = 6268.255
$ws.Range("I136").Value = 3324
$ws.Range("J136").Value = 10831.85
$ws.Range("K136").Value = 9972
$ws.Range("L136").Value = 32495.55
$ws.Range("M136").Value = -7422
$ws.Range("N136").Value = -37595.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3041.5789
$ws.Range("I3").Value = 2893
$ws.Range("K3").Value = 2893
$ws.Range("M3").Value = -2779
$ws.Range("H86").Value = 125006184
$ws.Range("I86").Value = 5365.75
$ws.Range("K86").Value = 5365.75
$ws.Range("M86").Value = -4242.75
$ws.Range("H89").Value = 125006184
$ws.Range("I89").Value = 5365.75
$ws.Range("K89").Value = 26828.75
$ws.Range("M89").Value = -21212.75
$ws.Range("H94").Value = 19232876
$ws.Range("I94").Value = 26316734
$ws.Range("K94").Value = 26316734
$ws.Range("M94").Value = -26316283
$ws.Range("H99").Value = 2029.1786
$ws.Range("I99").Value = 2191.625
$ws.Range("J99").Value = 1054.5
$ws.Range("K99").Value = 2191.625
$ws.Range("L99").Value = 1054.5
$ws.Range("M99").Value = -693.625
$ws.Range("N99").Value = -4050.5
$ws.Range("H105").Value = 4480.1333
$ws.Range("I105").Value = 2872.5
$ws.Range("K105").Value = 2872.5
$ws.Range("M105").Value = -1125.5
$ws.Range("H134").Value = 4209.0835
$ws.Range("I134").Value = 888.74194
$ws.Range("K134").Value = 2666.22582
$ws.Range("M134").Value = -131.2258200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4741.5312
$ws.Range("I16").Value = 4091.25
$ws.Range("K16").Value = 4091.25
$ws.Range("M16").Value = -3804.25
$ws.Range("H22").Value = 376.7857
$ws.Range("I22").Value = 289.58334
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 289.58334
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 60.41665999999998
$ws.Range("N22").Value = -1600
$ws.Range("H31").Value = 9725.226000000001
$ws.Range("I31").Value = 5230
$ws.Range("J31").Value = 11865.81
$ws.Range("K31").Value = 5230
$ws.Range("L31").Value = 11865.81
$ws.Range("M31").Value = -4935
$ws.Range("N31").Value = -12455.81
$ws.Range("H34").Value = 9725.226000000001
$ws.Range("I34").Value = 5230
$ws.Range("J34").Value = 11865.81
$ws.Range("K34").Value = 5230
$ws.Range("L34").Value = 11865.81
$ws.Range("M34").Value = -5028
$ws.Range("N34").Value = -12269.81
$ws.Range("H113").Value = 4741.5312
$ws.Range("I113").Value = 4091.25
$ws.Range("K113").Value = 4091.25
$ws.Range("M113").Value = -1921.25
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 16669264
$ws.Range("I14").Value = 16669264
$ws.Range("K14").Value = 50007792
$ws.Range("M14").Value = -50007619
$ws.Range("H40").Value = 123.3
$ws.Range("I40").Value = 75.75
$ws.Range("K40").Value = 303
$ws.Range("M40").Value = -234
$ws.Range("H139").Value = 107401.8
$ws.Range("I139").Value = 145221.58
$ws.Range("K139").Value = 435664.74
$ws.Range("M139").Value = -430524.74

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 38504480
$ws.Range("I122").Value = 52688036
$ws.Range("K122").Value = 158064108
$ws.Range("M122").Value = -158061658
$ws.Range("H126").Value = 6292.2856
$ws.Range("I126").Value = 4478.4
$ws.Range("K126").Value = 13435.2
$ws.Range("M126").Value = -10965.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5546.923
$ws.Range("I40").Value = 5357
$ws.Range("K40").Value = 5357
$ws.Range("M40").Value = -5221
$ws.Range("H46").Value = 4659.6665
$ws.Range("J46").Value = 4102.3335
$ws.Range("L46").Value = 4102.3335
$ws.Range("N46").Value = -4478.3335
$ws.Range("H132").Value = 15631527
$ws.Range("I132").Value = 31252748
$ws.Range("J132").Value = 10305.8125
$ws.Range("K132").Value = 93758244
$ws.Range("L132").Value = 30917.4375
$ws.Range("M132").Value = -93755714
$ws.Range("N132").Value = -35977.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 678.52
$ws.Range("I107").Value = 563.4737
$ws.Range("J107").Value = 1042.8334
$ws.Range("K107").Value = 1690.4211
$ws.Range("L107").Value = 3128.5002
$ws.Range("M107").Value = 229.5789
$ws.Range("N107").Value = -6968.5002
$ws.Range("H122").Value = 4337.933
$ws.Range("I122").Value = 3420.7368
$ws.Range("J122").Value = 5922.1816
$ws.Range("K122").Value = 10262.2104
$ws.Range("L122").Value = 17766.5448
$ws.Range("M122").Value = -7812.2104
$ws.Range("N122").Value = -22666.5448
$ws.Range("H132").Value = 23813000
$ws.Range("I132").Value = 29414118
$ws.Range("J132").Value = 8249.75
$ws.Range("K132").Value = 88242354
$ws.Range("L132").Value = 24749.25
$ws.Range("M132").Value = -88239824
$ws.Range("N132").Value = -29809.25
